# optimise the buff immune system
#
# MonsterRace.xlsx table2 row for the "鸟" (Bird) race had its buff-immunity
# description and a few related counters tweaked:
#   - Des (column C) changes from "擅长快速移动/擅长冲锋技能" (fast movement/
#     charge) to "擅长克制技能/擅长冲锋技能" (counter/charge)
#   - CountRush (E), CountOver (J) and CountHeal (Q) counters are adjusted
#     to match the new description
# Two unrelated AoE-count corrections are also applied (精灵/Elf row 6 and
# 植物/Plant row 17).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9 - 鸟 (Bird): update description + related skill counters
$ws.Range("C9").Value = "--擅长克制技能`$--擅长冲锋技能"
$ws.Range("E9").Value = 1
$ws.Range("J9").Value = 5
$ws.Range("Q9").Value = 1

# Row 6 - 精灵 (Elf): CountAoe correction
$ws.Range("H6").Value = 1

# Row 17 - 植物 (Plant): CountAoe correction
$ws.Range("H17").Value = 1
